$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new row 18 ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A18").Value = "Uitnodiging voor netwerkevent"
$ws.Range("B18").Value = "mailmind.test@zohomail.eu"
$ws.Range("C18").Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$ws.Range("D18").Value = "Samenwerking / Partnerverzoek"
$ws.Range("F18").Value = "2025-06-19 21:35:10"
$ws.Range("G18").Value = "Nee"

# Extend the conditional-formatting ranges to cover the new row
$dCond = $ws.Range("D2:D17").FormatConditions
for ($i = 1; $i -le $dCond.Count; $i++) {
    $dCond.Item($i).ModifyAppliesToRange($ws.Range("D2:D18"))
}

$gCond = $ws.Range("G2:G17").FormatConditions
for ($i = 1; $i -le $gCond.Count; $i++) {
    $gCond.Item($i).ModifyAppliesToRange($ws.Range("G2:G18"))
}

# --- "Dashboard" sheet: bump the "Samenwerking / Partnerverzoek" count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 5

Write-Output "done"
